$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (e.g. "13.80", "66.009.88") and must
# remain plain text, matching the source inlineStr cells. Force text format before
# assigning so Excel does not auto-convert them to numbers, then restore the
# original (unstyled) cell style so no formatting diff is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.932.76'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.446.62'
$ws.Range("D3").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.46'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.82'
$ws.Range("D6").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.443.86'
$ws.Range("D9").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.041.25'
$ws.Range("D13").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.88'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.062.35'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000170'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.444.81'
$ws.Range("D18").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.80'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.94'
$ws.Range("D21").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.91'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.531'
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.64'
$ws.Range("D27").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.13'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '29.11'
$ws.Range("D38").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.741.88'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.55'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.38'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.29'
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.98'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.21'
$ws.Range("D47").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '304.40'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.812'
$ws.Range("D51").Style = "Normal"

# Columns B, C and E are plain text already (names, links, padded percentages)
# and do not get coerced to numbers, so a direct Value assignment is sufficient.
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("E6").Value = '  -4.79%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -4.47%  '
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("E10").Value = '  -6.96%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("E12").Value = '  -4.68%  '
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("E15").Value = '  -6.66%  '
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("E17").Value = '  -3.92%  '
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("E19").Value = '  -5.07%  '
$ws.Range("E20").Value = '  -2.18%  '
$ws.Range("E21").Value = '  -7.23%  '
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.60%  '
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  -7.30%  '
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("E32").Value = '  -6.67%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  -8.48%  '
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("E37").Value = '  -1.73%  '
$ws.Range("E38").Value = '  +10.85%  '
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("E40").Value = '  -4.82%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E42").Value = '  -10.88%  '
$ws.Range("E43").Value = '  -6.16%  '
$ws.Range("E44").Value = '  -7.25%  '
$ws.Range("E45").Value = '  -4.86%  '
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("E47").Value = '  -7.86%  '
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E50").Value = '  -3.78%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("E51").Value = '  -4.12%  '
